$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 212
$ws.Cells.Item(212, 1).Value = "(Intercept)"
$ws.Cells.Item(212, 2).Value = [double]"2.739275404261982e-112"
$ws.Cells.Item(212, 3).Value = [double]"24.77271601144333"
$ws.Cells.Item(212, 4).Value = [double]"-10.36954675709533"
$ws.Cells.Item(212, 5).Value = [double]"3.411389371628055e-25"
$ws.Cells.Item(212, 6).Value = [double]"2.024051504041591e-133"
$ws.Cells.Item(212, 7).Value = [double]"3.060408745163685e-91"
$ws.Cells.Item(212, 8).Value = "TZP"

# Row 213
$ws.Cells.Item(213, 1).Value = "Year"
$ws.Cells.Item(213, 2).Value = [double]"1.13609995950807"
$ws.Cells.Item(213, 3).Value = [double]"0.01227980250051377"
$ws.Cells.Item(213, 4).Value = [double]"10.39115319218002"
$ws.Cells.Item(213, 5).Value = [double]"2.720437656928805e-25"
$ws.Cells.Item(213, 6).Value = [double]"1.109131770857509"
$ws.Cells.Item(213, 7).Value = [double]"1.163835233989738"
$ws.Cells.Item(213, 8).Value = "TZP"

# Row 214
$ws.Cells.Item(214, 1).Value = "Specimen_typeother"
$ws.Cells.Item(214, 2).Value = [double]"0.7873899537786169"
$ws.Cells.Item(214, 3).Value = [double]"0.1110123134232805"
$ws.Cells.Item(214, 4).Value = [double]"-2.153199513725651"
$ws.Cells.Item(214, 5).Value = [double]"0.03130300188647972"
$ws.Cells.Item(214, 6).Value = [double]"0.6331707152625055"
$ws.Cells.Item(214, 7).Value = [double]"0.9784872850853259"
$ws.Cells.Item(214, 8).Value = "TZP"

# Row 215
$ws.Cells.Item(215, 1).Value = "Specimen_typeRespiratory"
$ws.Cells.Item(215, 2).Value = [double]"0.5975810744520984"
$ws.Cells.Item(215, 3).Value = [double]"0.1195533034007147"
$ws.Cells.Item(215, 4).Value = [double]"-4.306575395881504"
$ws.Cells.Item(215, 5).Value = [double]"1.658014335821558e-05"
$ws.Cells.Item(215, 6).Value = [double]"0.4723698948825547"
$ws.Cells.Item(215, 7).Value = [double]"0.7548688973041973"
$ws.Cells.Item(215, 8).Value = "TZP"

# Row 216
$ws.Cells.Item(216, 1).Value = "Specimen_typeUrine"
$ws.Cells.Item(216, 2).Value = [double]"0.6283244116518715"
$ws.Cells.Item(216, 3).Value = [double]"0.09542346174286635"
$ws.Cells.Item(216, 4).Value = [double]"-4.869857562741592"
$ws.Cells.Item(216, 5).Value = [double]"1.116787225151061e-06"
$ws.Cells.Item(216, 6).Value = [double]"0.5212768242716972"
$ws.Cells.Item(216, 7).Value = [double]"0.7577947528162462"
$ws.Cells.Item(216, 8).Value = "TZP"

# Row 217
$ws.Cells.Item(217, 1).Value = "Specimen_typeWound & soft tissues"
$ws.Cells.Item(217, 2).Value = [double]"0.8616304390094731"
$ws.Cells.Item(217, 3).Value = [double]"0.1036174040281501"
$ws.Cells.Item(217, 4).Value = [double]"-1.437295468759839"
$ws.Cells.Item(217, 5).Value = [double]"0.1506340558382742"
$ws.Cells.Item(217, 6).Value = [double]"0.7031611904248042"
$ws.Cells.Item(217, 7).Value = [double]"1.055577748291711"
$ws.Cells.Item(217, 8).Value = "TZP"

# Row 218
$ws.Cells.Item(218, 1).Value = "HospitalCHBH"
$ws.Cells.Item(218, 2).Value = [double]"0.2469828767240796"
$ws.Cells.Item(218, 3).Value = [double]"0.2186812637033222"
$ws.Cells.Item(218, 4).Value = [double]"-6.394860931744881"
$ws.Cells.Item(218, 5).Value = [double]"1.606935406072353e-10"
$ws.Cells.Item(218, 6).Value = [double]"0.1600109447956319"
$ws.Cells.Item(218, 7).Value = [double]"0.3774263422482198"
$ws.Cells.Item(218, 8).Value = "TZP"

# Row 219
$ws.Cells.Item(219, 1).Value = "HospitalCNGMO"
$ws.Cells.Item(219, 2).Value = [double]"0.8551143957151932"
$ws.Cells.Item(219, 3).Value = [double]"0.4581951089018857"
$ws.Cells.Item(219, 4).Value = [double]"-0.3416012519476278"
$ws.Cells.Item(219, 5).Value = [double]"0.7326509949070338"
$ws.Cells.Item(219, 6).Value = [double]"0.3414622350443337"
$ws.Cells.Item(219, 7).Value = [double]"2.100357510593966"
$ws.Cells.Item(219, 8).Value = "TZP"

# Row 220
$ws.Cells.Item(220, 1).Value = "HospitalTCB"
$ws.Cells.Item(220, 2).Value = [double]"0.7374625541100627"
$ws.Cells.Item(220, 3).Value = [double]"0.1635997794239589"
$ws.Cells.Item(220, 4).Value = [double]"-1.861493745392306"
$ws.Cells.Item(220, 5).Value = [double]"0.06267448307947809"
$ws.Cells.Item(220, 6).Value = [double]"0.5343568961482102"
$ws.Cells.Item(220, 7).Value = [double]"1.015326189401102"
$ws.Cells.Item(220, 8).Value = "TZP"

# Row 221
$ws.Cells.Item(221, 1).Value = "Ward_ED_ICUED"
$ws.Cells.Item(221, 2).Value = [double]"0.1751083169337728"
$ws.Cells.Item(221, 3).Value = [double]"0.2294954369635499"
$ws.Cells.Item(221, 4).Value = [double]"-7.592092312062168"
$ws.Cells.Item(221, 5).Value = [double]"3.147806244090618e-14"
$ws.Cells.Item(221, 6).Value = [double]"0.1108873571682885"
$ws.Cells.Item(221, 7).Value = [double]"0.2729098799616956"
$ws.Cells.Item(221, 8).Value = "TZP"

# Row 222
$ws.Cells.Item(222, 1).Value = "Ward_ED_ICUOther"
$ws.Cells.Item(222, 2).Value = [double]"0.2404926322365203"
$ws.Cells.Item(222, 3).Value = [double]"0.1611208584379195"
$ws.Cells.Item(222, 4).Value = [double]"-8.844701045618217"
$ws.Cells.Item(222, 5).Value = [double]"9.177367799809633e-19"
$ws.Cells.Item(222, 6).Value = [double]"0.1750642547097513"
$ws.Cells.Item(222, 7).Value = [double]"0.3294294254595298"
$ws.Cells.Item(222, 8).Value = "TZP"

# Row 223
$ws.Cells.Item(223, 1).Value = "GenderF"
$ws.Cells.Item(223, 2).Value = [double]"0.8291972810592688"
$ws.Cells.Item(223, 3).Value = [double]"0.06364523514550369"
$ws.Cells.Item(223, 4).Value = [double]"-2.942831101030215"
$ws.Cells.Item(223, 5).Value = [double]"0.00325225799907553"
$ws.Cells.Item(223, 6).Value = [double]"0.7319187350048831"
$ws.Cells.Item(223, 7).Value = [double]"0.939349477382989"
$ws.Cells.Item(223, 8).Value = "TZP"

# Row 224
$ws.Cells.Item(224, 1).Value = "HospitalCHBH:Ward_ED_ICUED"
$ws.Cells.Item(224, 2).Value = [double]"3.042335421895682"
$ws.Cells.Item(224, 3).Value = [double]"0.4159323903694834"
$ws.Cells.Item(224, 4).Value = [double]"2.675015163830329"
$ws.Cells.Item(224, 5).Value = [double]"0.007472587019224445"
$ws.Cells.Item(224, 6).Value = [double]"1.311135616803233"
$ws.Cells.Item(224, 7).Value = [double]"6.750207097466561"
$ws.Cells.Item(224, 8).Value = "TZP"

# Row 225
$ws.Cells.Item(225, 1).Value = "HospitalCNGMO:Ward_ED_ICUED"
$ws.Cells.Item(225, 8).Value = "TZP"

# Row 226
$ws.Cells.Item(226, 1).Value = "HospitalTCB:Ward_ED_ICUED"
$ws.Cells.Item(226, 2).Value = [double]"1.415772429663642"
$ws.Cells.Item(226, 3).Value = [double]"0.5265116355960461"
$ws.Cells.Item(226, 4).Value = [double]"0.6603372942513585"
$ws.Cells.Item(226, 5).Value = [double]"0.5090374024045634"
$ws.Cells.Item(226, 6).Value = [double]"0.4686435918126847"
$ws.Cells.Item(226, 7).Value = [double]"3.799266388146391"
$ws.Cells.Item(226, 8).Value = "TZP"

# Row 227
$ws.Cells.Item(227, 1).Value = "HospitalCHBH:Ward_ED_ICUOther"
$ws.Cells.Item(227, 2).Value = [double]"3.472528998949473"
$ws.Cells.Item(227, 3).Value = [double]"0.242570149906032"
$ws.Cells.Item(227, 4).Value = [double]"5.13205415886129"
$ws.Cells.Item(227, 5).Value = [double]"2.865970000453856e-07"
$ws.Cells.Item(227, 6).Value = [double]"2.165756455575788"
$ws.Cells.Item(227, 7).Value = [double]"5.609057043838499"
$ws.Cells.Item(227, 8).Value = "TZP"

# Row 228
$ws.Cells.Item(228, 1).Value = "HospitalCNGMO:Ward_ED_ICUOther"
$ws.Cells.Item(228, 2).Value = [double]"3.573601659993845"
$ws.Cells.Item(228, 3).Value = [double]"0.5006997089511203"
$ws.Cells.Item(228, 4).Value = [double]"2.543588368044857"
$ws.Cells.Item(228, 5).Value = [double]"0.01097202997634257"
$ws.Cells.Item(228, 6).Value = [double]"1.33916432239849"
$ws.Cells.Item(228, 7).Value = [double]"9.686798314095258"
$ws.Cells.Item(228, 8).Value = "TZP"

# Row 229
$ws.Cells.Item(229, 1).Value = "HospitalTCB:Ward_ED_ICUOther"
$ws.Cells.Item(229, 2).Value = [double]"3.467871074136891"
$ws.Cells.Item(229, 3).Value = [double]"0.1852223581114965"
$ws.Cells.Item(229, 4).Value = [double]"6.713773083615963"
$ws.Cells.Item(229, 5).Value = [double]"1.896551695123008e-11"
$ws.Cells.Item(229, 6).Value = [double]"2.413805353123875"
$ws.Cells.Item(229, 7).Value = [double]"4.991314100864938"
$ws.Cells.Item(229, 8).Value = "TZP"
